# Fruta / hortaliza, semanal
# Insert a new weekly record as row 4 (pushing all subsequent records down
# by one row) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4..46 down to 5..47, duplicating row 3's formatting into the
# freshly-opened row 4 (matches Excel's native Insert-row behaviour).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with this week's observation.
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44691
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 100112030
$ws.Range("G4").Value = "Poroto granado"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 25000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 25000
$ws.Range("N4").Value = "`$/saco 25 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 1000
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
